# Append 5 new backtest rows (58-62) to the prediction log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 58; A = "2023-03-22-09:00"; B = 28071.6;  C = 28120.7;  D = 27962.1;  E = 28086.7;  F = 28003.57721159537;  G = 27999.57147313305;  H = 18862.109 },
    @{ Row = 59; A = "2023-03-22-10:00"; B = 28086.7;  C = 28178.6;  D = 28049.1;  E = 28138.8;  F = 27999.57147313305;  G = 28124.06128157466;  H = 13991.948 },
    @{ Row = 60; A = "2023-03-22-11:00"; B = 28138.8;  C = 28197.9;  D = 28101.8;  E = 28113.5;  F = 28124.06128157466;  G = 28099.63101356567;  H = 10376.928 },
    @{ Row = 61; A = "2023-03-22-12:00"; B = 28113.5;  C = 28241;    D = 28090.4;  E = 28207.1;  F = 28099.63101356567;  G = 28283.43265622613;  H = 19141.998 },
    @{ Row = 62; A = "2023-03-22-13:00"; B = 28207;    C = 28390;    D = 28150.5;  E = 28269.7;  F = 28283.43265622613;  G = 28399.89992704705;  H = 38841.696 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
}
